$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos list values (prices / volume %) scraped on
# Thu Sep 21 23:20:10 UTC 2023. Leading apostrophe forces text entry so
# purely-numeric-looking price strings (e.g. "211.19") are not coerced
# into numeric cells; Style reset avoids leaving a stray @ text format.

$ws.Range("D2").Value = "'26.625.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.82%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.588.38"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.07%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.05%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'211.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.23%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'0.509"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -2.34%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.06%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.246"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.02%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.62%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'19.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -3.37%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.40%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.811.28"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.05%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.583.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -2.38%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  -2.63%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  -3.76%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'64.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.27%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'26.611.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.75%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.0₃0728"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.98%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'209.15"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.11%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +0.04%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.97%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.47%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'2.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.60%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'8.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.96%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'146.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.32%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +0.00%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -0.62%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -2.65%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'15.30"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.73%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'0.0507"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +0.23%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -1.39%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.60%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.677"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +22.36%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'2.90"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.01%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.311.96"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.04%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.49%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -5.01%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -2.92%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.827"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.58%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +0.05%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'TrustWalletToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.791"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.36%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'FraxShare"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'5.38"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +2.89%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -2.70%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'62.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.43%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.724.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -1.86%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'89.48"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.99%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.66%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.841"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -8.04%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  -1.59%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0979"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -1.55%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'7.55"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.22%  "
$ws.Range("E51").Style = "Normal"
